$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Extend the existing truth-table rows 18-24 with three more input
# columns (G:K) that were added to the table.
# ---------------------------------------------------------------------
$ghijk = @(
    @(0,0,0,0,0),
    @(0,0,0,1,0),
    @(0,0,1,0,0),
    @(0,0,1,1,0),
    @(1,1,1,1,0),
    @(1,1,1,0,0),
    @(1,1,0,1,0)
)

$rows = $ghijk.Length
$cols = 5
$arr = New-Object 'object[,]' $rows,$cols
for ($i = 0; $i -lt $rows; $i++) {
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i,$j] = $ghijk[$i][$j]
    }
}
$ws.Range("G18:K24").Value = $arr

# ---------------------------------------------------------------------
# Add a new truth-table block in rows 26-32 (columns D:F).
# ---------------------------------------------------------------------
$def = @(
    @(0,0,0),
    @(0,0,1),
    @(0,1,0),
    @(0,1,1),
    @(1,1,1),
    @(1,1,0),
    @(1,0,1)
)

$rows2 = $def.Length
$cols2 = 3
$arr2 = New-Object 'object[,]' $rows2,$cols2
for ($i = 0; $i -lt $rows2; $i++) {
    for ($j = 0; $j -lt $cols2; $j++) {
        $arr2[$i,$j] = $def[$i][$j]
    }
}
$ws.Range("D26:F32").Value = $arr2

# ---------------------------------------------------------------------
# Update the worksheet view: scroll back up and select the newly-added
# block, matching the saved view state in the workbook.
# ---------------------------------------------------------------------
$ws.Range("D26:F32").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
